$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(4.9000000000000004, 10, 11, 15, 16, 1, 1, 1, 1, 0.04, 0.04, 0.04, 0.04, 200, 300, 600, 0.15049296094977052, 0.0013721046507395801, 0.3480414117524081, 10000)
$row3 = @(4.9000000000000004, 10, 11, 15, 16, 1, 1, 1, 1, 0.04, 0.04, 0.04, 0.04, 200, 300, 600, 0.19161495366676654, 0.0012558213135784916, 0.33438156265954805, 10000)

for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$ws.Range("A3:T3").Select()
